$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Il27 -> Il27ra ligand-receptor pair data (new TPM values),
# now including "Inflammatory-Mac" as both a sending and target cluster,
# and no longer a "Resolving-Mac" -> "Resolving-Mac" self-pair row.
$data = @(
    @("FAPs", "Il27", "Il27ra", "ECs", 3, 1, 0.9134243333333334, 2.740273, 0.2338231012639866, 0.2338231012639866, 3, 1, 2.171010666666667, 6.513032, 0.5189085589295733, 0.5189085589295733, 1.983053970859556, 17.847485737736, 0.121332808521339, 0.121332808521339),
    @("FAPs", "Il27", "Il27ra", "FAPs", 3, 1, 0.9134243333333334, 2.740273, 0.2338231012639866, 0.2338231012639866, 3, 1, 1.561521333333333, 4.684564, 0.3732302181922886, 0.3732302181922886, 1.426331582885778, 12.836984245972, 0.0872698471031553, 0.08726984710315532),
    @("FAPs", "Il27", "Il27ra", "Inflammatory-Mac", 3, 1, 0.9134243333333334, 2.740273, 0.2338231012639866, 0.2338231012639866, 3, 1, 0.4199566666666667, 1.25987, 0.1003768024076347, 0.1003768024076346, 0.383598638278889, 3.45238774451, 0.02347041523391553, 0.02347041523391553),
    @("FAPs", "Il27", "Il27ra", "MuSCs", 3, 1, 0.9134243333333334, 2.740273, 0.2338231012639866, 0.2338231012639866, 1, 0.3333333333333333, 0.03131333333333333, 0.09394, 0.007484420470503464, 0.007484420470503464, 0.02860236062444444, 0.25742124562, 0.001750030405576786, 0.001750030405576786),
    @("Inflammatory-Mac", "Il27", "Il27ra", "ECs", 3, 1, 1.830769666666667, 5.492309000000001, 0.4686499204568687, 0.4686499204568688, 3, 1, 2.171010666666667, 6.513032, 0.5189085589295733, 0.5189085589295733, 3.974620474543111, 35.77158427088801, 0.2431864548667329, 0.2431864548667329),
    @("Inflammatory-Mac", "Il27", "Il27ra", "FAPs", 3, 1, 1.830769666666667, 5.492309000000001, 0.4686499204568687, 0.4686499204568688, 3, 1, 1.561521333333333, 4.684564, 0.3732302181922886, 0.3732302181922886, 2.858785890919556, 25.729073018276, 0.1749143120679158, 0.1749143120679158),
    @("Inflammatory-Mac", "Il27", "Il27ra", "Inflammatory-Mac", 3, 1, 1.830769666666667, 5.492309000000001, 0.4686499204568687, 0.4686499204568688, 3, 1, 0.4199566666666667, 1.25987, 0.1003768024076347, 0.1003768024076346, 0.7688439266477779, 6.919595339830001, 0.04704158046405281, 0.04704158046405281),
    @("Inflammatory-Mac", "Il27", "Il27ra", "MuSCs", 3, 1, 1.830769666666667, 5.492309000000001, 0.4686499204568687, 0.4686499204568688, 1, 0.3333333333333333, 0.03131333333333333, 0.09394, 0.007484420470503464, 0.007484420470503464, 0.05732750082888889, 0.5159475074600001, 0.003507573058167208, 0.003507573058167208),
    @("Resolving-Mac", "Il27", "Il27ra", "ECs", 3, 1, 1.162282, 3.486846, 0.2975269782791446, 0.2975269782791446, 3, 1, 2.171010666666667, 6.513032, 0.5189085589295733, 0.5189085589295733, 2.523326619674667, 22.709939577072, 0.1543892955415014, 0.1543892955415014),
    @("Resolving-Mac", "Il27", "Il27ra", "FAPs", 3, 1, 1.162282, 3.486846, 0.2975269782791446, 0.2975269782791446, 3, 1, 1.561521333333333, 4.684564, 0.3732302181922886, 0.3732302181922886, 1.814928138349333, 16.334353245144, 0.1110460590212175, 0.1110460590212175),
    @("Resolving-Mac", "Il27", "Il27ra", "Inflammatory-Mac", 3, 1, 1.162282, 3.486846, 0.2975269782791446, 0.2975269782791446, 3, 1, 0.4199566666666667, 1.25987, 0.1003768024076347, 0.1003768024076346, 0.4881080744466667, 4.39297267002, 0.02986480670966631, 0.02986480670966631),
    @("Resolving-Mac", "Il27", "Il27ra", "MuSCs", 3, 1, 1.162282, 3.486846, 0.2975269782791446, 0.2975269782791446, 1, 0.3333333333333333, 0.03131333333333333, 0.09394, 0.007484420470503464, 0.007484420470503464, 0.03639492369333333, 0.32755431324, 0.00222681700675947, 0.00222681700675947)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowData[$j]
    }
}